# Updates the cryptos price list (Price / Volume(1h) columns) with refreshed
# market data, and swaps the Aptos / RenzoRestakedETH rows (36 and 37).
#
# All Price/Coin/Link values in this sheet are stored as literal text
# (inlineStr) rather than numbers. Several of the new Price values look like
# plain numbers (e.g. "7.59"), and Excel would normally auto-convert such a
# value to a numeric cell when it is assigned. To keep those cells as text
# (matching the original file's cell typing) we temporarily mark them with a
# Text ("@") number format before writing the value, then restore the
# default "Normal" style afterwards so no visible/applied formatting change
# is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-recognized as a number by Excel.
$textCells = @("D4","D5","D6","D9","D10","D12","D14","D19","D20","D21","D22","D23","D24","D25","D32","D35","D37","D40","D41","D42","D43","D46","D47","D50","D51")

foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "61.338.32"
$ws.Range("E2").Value = "  -0.58%  "
# Row 3 - Ethereum
$ws.Range("D3").Value = "3.438.78"
$ws.Range("E3").Value = "  -0.22%  "
# Row 4 - TetherUSD
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.12%  "
# Row 5 - BNB
$ws.Range("D5").Value = "574.00"
$ws.Range("E5").Value = "  -0.58%  "
# Row 6 - Solana
$ws.Range("D6").Value = "144.19"
$ws.Range("E6").Value = "  -2.60%  "
# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.438.34"
$ws.Range("E7").Value = "  -0.22%  "
# Row 9 - XRP
$ws.Range("D9").Value = "0.479"
$ws.Range("E9").Value = "  +0.92%  "
# Row 10 - Toncoin
$ws.Range("D10").Value = "7.59"
$ws.Range("E10").Value = "  -1.44%  "
# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +2.02%  "
# Row 12 - Cardano
$ws.Range("D12").Value = "0.389"
$ws.Range("E12").Value = "  -0.97%  "
# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "4.026.86"
$ws.Range("E13").Value = "  -0.06%  "
# Row 14 - Avalanche
$ws.Range("D14").Value = "28.49"
$ws.Range("E14").Value = "  +4.70%  "
# Row 15 - TRON
$ws.Range("E15").Value = "  -0.50%  "
# Row 16 - ShibaInu
$ws.Range("E16").Value = "  -1.94%  "
# Row 17 - WrappedEther
$ws.Range("D17").Value = "3.437.05"
$ws.Range("E17").Value = "  -0.24%  "
# Row 18 - WrappedBTC
$ws.Range("D18").Value = "61.456.64"
$ws.Range("E18").Value = "  -0.42%  "
# Row 19 - Polkadot
$ws.Range("D19").Value = "6.35"
$ws.Range("E19").Value = "  +4.40%  "
# Row 20 - Chainlink
$ws.Range("D20").Value = "14.37"
$ws.Range("E20").Value = "  +1.84%  "
# Row 21 - Uniswap
$ws.Range("D21").Value = "9.34"
$ws.Range("E21").Value = "  -2.66%  "
# Row 22 - BitcoinCash
$ws.Range("D22").Value = "400.01"
$ws.Range("E22").Value = "  +3.77%  "
# Row 23 - Polygon
$ws.Range("D23").Value = "0.563"
$ws.Range("E23").Value = "  +0.23%  "
# Row 24 - Litecoin
$ws.Range("D24").Value = "74.07"
$ws.Range("E24").Value = "  +2.98%  "
# Row 25 - Dai
$ws.Range("D25").Value = "0.998"
$ws.Range("E25").Value = "  -0.45%  "
# Row 26 - PEPE
$ws.Range("E26").Value = "  -3.32%  "
# Row 27 - WrappedeETH
$ws.Range("D27").Value = "3.583.40"
$ws.Range("E27").Value = "  +0.49%  "
# Row 28 - Kaspa
$ws.Range("E28").Value = "  +0.56%  "
# Row 29 - RenderToken
$ws.Range("E29").Value = "  -2.77%  "
# Row 30 - Binance-PegBSC-USD
$ws.Range("E30").Value = "  +0.10%  "
# Row 31 - Fetch.AI
$ws.Range("E31").Value = "  -6.77%  "
# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").Value = "8.21"
$ws.Range("E32").Value = "  -0.24%  "
# Row 33 - PancakeSwap
$ws.Range("E33").Value = "  +0.19%  "
# Row 34 - USDe
$ws.Range("E34").Value = "  +0.00%  "
# Row 35 - EthereumClassic
$ws.Range("D35").Value = "23.97"
$ws.Range("E35").Value = "  +0.34%  "

# Rows 36/37 - Aptos and RenzoRestakedETH swap places (with refreshed data)
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").Value = "3.468.16"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").Value = "7.02"
$ws.Range("E37").Value = "  -0.37%  "

# Row 38 - NEARProtocol
$ws.Range("E38").Value = "  -3.45%  "
# Row 39 - ImmutableX
$ws.Range("E39").Value = "  -1.85%  "
# Row 40 - Monero
$ws.Range("D40").Value = "166.65"
$ws.Range("E40").Value = "  +0.08%  "
# Row 41 - Hedera
$ws.Range("D41").Value = "0.0789"
$ws.Range("E41").Value = "  -0.66%  "
# Row 42 - EnergySwap
$ws.Range("D42").Value = "27.41"
$ws.Range("E42").Value = "  +4.19%  "
# Row 43 - Mantle
$ws.Range("D43").Value = "0.803"
$ws.Range("E43").Value = "  +1.61%  "
# Row 44 - Filecoin
$ws.Range("E44").Value = "  +1.35%  "
# Row 45 - FirstDigitalUSD
$ws.Range("E45").Value = "  +0.08%  "
# Row 46 - OKB
$ws.Range("D46").Value = "42.20"
$ws.Range("E46").Value = "  +0.02%  "
# Row 47 - Stacks
$ws.Range("D47").Value = "1.72"
$ws.Range("E47").Value = "  -1.38%  "
# Row 48 - Maker
$ws.Range("D48").Value = "2.616.47"
$ws.Range("E48").Value = "  -1.60%  "
# Row 49 - ONDO
$ws.Range("E49").Value = "  -5.67%  "
# Row 50 - Cosmos
$ws.Range("D50").Value = "6.95"
$ws.Range("E50").Value = "  +1.20%  "
# Row 51 - InjectiveProtocol
$ws.Range("D51").Value = "23.20"
$ws.Range("E51").Value = "  -3.66%  "

# Restore default styling on cells we temporarily forced to Text format.
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = "Normal"
}
